# Auto-generated edit script for "Add data for 2024-06-28"
# Updates year-2024 (and a few year-2021 correction) cell values
# across the Citywide Totals, By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3805
$ws.Range("K3").Value = 3848
$ws.Range("H4").Value = 1735
$ws.Range("K4").Value = 782
$ws.Range("K6").Value = 4386
$ws.Range("H7").Value = 26048
$ws.Range("K7").Value = 13087

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 255
$ws.Range("K3").Value = 269
$ws.Range("K6").Value = 291
$ws.Range("K7").Value = 885

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 279

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 146
$ws.Range("K7").Value = 539

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 150
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 440

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K6").Value = 123
$ws.Range("K7").Value = 304

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 63
$ws.Range("K3").Value = 94
$ws.Range("K7").Value = 232

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K5").Value = 2
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 110
$ws.Range("K4").Value = 45
$ws.Range("K5").Value = 27
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 381
$ws.Range("K8").Value = 885
$ws.Range("K11").Value = 265
$ws.Range("K12").Value = 23
$ws.Range("K18").Value = 89
$ws.Range("K19").Value = 405
$ws.Range("K20").Value = 294
$ws.Range("K24").Value = 42
$ws.Range("K27").Value = 130
$ws.Range("K29").Value = 693
$ws.Range("K30").Value = 47
$ws.Range("K31").Value = 143
$ws.Range("K33").Value = 539
$ws.Range("K37").Value = 440
$ws.Range("K41").Value = 111
$ws.Range("K42").Value = 460
$ws.Range("K44").Value = 122
$ws.Range("K47").Value = 76
$ws.Range("K48").Value = 167
$ws.Range("K49").Value = 73
$ws.Range("K51").Value = 150
$ws.Range("K54").Value = 248
$ws.Range("K55").Value = 148
$ws.Range("K58").Value = 7
$ws.Range("H63").Value = 287
$ws.Range("K64").Value = 79
$ws.Range("K65").Value = 304
$ws.Range("K67").Value = 514
$ws.Range("K72").Value = 62
$ws.Range("K73").Value = 118
$ws.Range("K78").Value = 161
$ws.Range("K79").Value = 339
$ws.Range("K80").Value = 47
$ws.Range("K83").Value = 279
$ws.Range("K84").Value = 95
$ws.Range("K85").Value = 594
$ws.Range("K89").Value = 180
$ws.Range("K94").Value = 161
$ws.Range("K97").Value = 110
$ws.Range("K98").Value = 68
$ws.Range("K99").Value = 232
$ws.Range("H101").Value = 26048
$ws.Range("K101").Value = 13087

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 154
$ws.Range("K3").Value = 175
$ws.Range("K7").Value = 514

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 248

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 197
$ws.Range("K3").Value = 242
$ws.Range("K6").Value = 200
$ws.Range("K7").Value = 693

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K3").Value = 36
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 131
$ws.Range("K3").Value = 120
$ws.Range("K7").Value = 405

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("K2").Value = 26
$ws.Range("K4").Value = 6
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 123
$ws.Range("K7").Value = 460

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("Dunning")
$ws.Range("K2").Value = 15
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 116
$ws.Range("K3").Value = 111
$ws.Range("K7").Value = 339

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 294

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 141
$ws.Range("K3").Value = 120
$ws.Range("K7").Value = 381

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K2").Value = 46
$ws.Range("K6").Value = 68
$ws.Range("K7").Value = 161

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 68

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 83
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 32
$ws.Range("K3").Value = 30
$ws.Range("K6").Value = 35
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 110

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 57
$ws.Range("K7").Value = 180

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 42
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 213
$ws.Range("K6").Value = 136
$ws.Range("K7").Value = 594

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 62

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K3").Value = 10
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 17
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 23

$ws = $wb.Worksheets.Item("Millenium Park")
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 7
